# The presentation currently carries two theme parts:
#   ppt/theme/theme1.xml  -> "Integral" color scheme (used by the slide master / all slides)
#   ppt/theme/theme2.xml  -> "Office Theme" color scheme (used by the notes master)
#
# The authored edit swaps the contents of these two theme parts: the slide
# master's theme becomes the "Office Theme" palette, and vice versa. Font
# scheme / format scheme are identical between the two themes, so the only
# observable difference is the 10 colors that differ (dk1/lt1 are already
# identical in both: 000000 / FFFFFF).
#
# Apply the "Office Theme" palette to the deck's (slide master) theme color
# scheme via the standard ThemeColorScheme COM object, in ThemeColorSchemeIndex
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      -> 000000
$cs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477    # folHlink -> 954F72
